$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old single header row (Force (N) / Voltage (V) / ADC)
$ws.Range("A1:C1").ClearContents()

# New layout:
#  Row 1: title cell in A1 only
#  Row 2: column headers in A2/B2
#  Rows 3-20: data pairs
# Order of assignment matters for shared-string index allocation, so write
# the header row cells before the title cell.
$ws.Range("A2").Value = "Load (kg)"
$ws.Range("B2").Value = "Voltage tst 1 (V)"
$ws.Range("A1").Value = "FSR_2 Input test"

$data = @(
    @(0,   3.46),
    @(0.1, 1.7),
    @(0.2, 0.56),
    @(0.3, 0.46),
    @(0.4, 0.35),
    @(0.5, 0.29),
    @(0.6, 0.28),
    @(0.7, 0.247),
    @(0.8, 0.21),
    @(0.9, 0.2),
    @(1,   0.195),
    @(1.5, 0.14),
    @(2,   0.11),
    @(2.5, 0.103),
    @(3,   0.101),
    @(3.5, 0.099),
    @(4,   0.079),
    @(4.5, 0.076)
)

$row = 3
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

# Column B/C width (bestFit widened to fit "Voltage tst 1 (V)")
$ws.Range("B:C").ColumnWidth = 13

# Best-effort: shrink/reposition the workbook window (matches the author's
# saved window state in the diff).
$win = $excel.ActiveWindow
$win.Left = -90
$win.Top = -90
$win.Width = 19380
$win.Height = 10260

# Selection moves to D3 in the saved file
$ws.Range("D3").Select() | Out-Null
